$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.982.27'
$ws.Range("E2").Value = '  -0.32%  '

$ws.Range("D3").Value = '2.573.58'
$ws.Range("E3").Value = '  +0.12%  '

$ws.Range("D4").Value2 = "'0.999"
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value2 = "'584.64"
$ws.Range("E5").Value = '  -0.05%  '

$ws.Range("D6").Value2 = "'144.22"
$ws.Range("E6").Value = '  -2.14%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  -1.95%  '

$ws.Range("E9").Value = '  -2.23%  '

$ws.Range("E10").Value = '  -0.83%  '

$ws.Range("E11").Value = '  -0.28%  '

$ws.Range("E12").Value = '  -1.99%  '

$ws.Range("D13").Value2 = "'26.98"
$ws.Range("E13").Value = '  -1.48%  '

$ws.Range("D14").Value = '3.032.62'
$ws.Range("E14").Value = '  +0.05%  '

$ws.Range("D15").Value = '62.886.02'
$ws.Range("E15").Value = '  -0.43%  '

$ws.Range("E16").Value = '  -1.76%  '

$ws.Range("D17").Value = '2.570.90'
$ws.Range("E17").Value = '  -1.41%  '

$ws.Range("D18").Value2 = "'11.03"
$ws.Range("E18").Value = '  -2.82%  '

$ws.Range("D19").Value2 = "'340.48"
$ws.Range("E19").Value = '  -0.71%  '

$ws.Range("E20").Value = '  -2.02%  '

$ws.Range("D21").Value2 = "'6.63"
$ws.Range("E21").Value = '  -3.65%  '

$ws.Range("E22").Value = '  +0.09%  '

$ws.Range("E23").Value = '  +3.42%  '

$ws.Range("D24").Value2 = "'67.65"
$ws.Range("E24").Value = '  +1.21%  '

$ws.Range("D25").Value2 = "'1.61"
$ws.Range("E25").Value = '  +7.70%  '

$ws.Range("D26").Value2 = "'1.61"
$ws.Range("E26").Value = '  -1.56%  '

$ws.Range("E27").Value = '  -3.71%  '

$ws.Range("D28").Value2 = "'7.99"
$ws.Range("E28").Value = '  -2.21%  '

$ws.Range("E29").Value = '  -0.23%  '

$ws.Range("D30").Value2 = "'8.24"
$ws.Range("E30").Value = '  -3.01%  '

$ws.Range("E31").Value = '  -2.77%  '

$ws.Range("D32").Value2 = "'463.57"
$ws.Range("E32").Value = '  -0.31%  '

$ws.Range("D33").Value = '0.0₃0797'
$ws.Range("E33").Value = '  -3.62%  '

$ws.Range("D34").Value2 = "'1.66"
$ws.Range("E34").Value = '  +1.63%  '

$ws.Range("D35").Value2 = "'176.36"
$ws.Range("E35").Value = '  +0.17%  '

$ws.Range("E36").Value = '  -0.04%  '

$ws.Range("E37").Value = '  -2.42%  '

$ws.Range("E38").Value = '  -2.18%  '

$ws.Range("D39").Value2 = "'4.55"
$ws.Range("E39").Value = '  +0.01%  '

$ws.Range("E40").Value = '  +0.05%  '

$ws.Range("D41").Value2 = "'1.69"

$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value2 = "'158.05"
$ws.Range("E42").Value = '  +4.29%  '

$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").Value2 = "'39.94"
$ws.Range("E43").Value = '  +0.94%  '

$ws.Range("D44").Value2 = "'3.68"
$ws.Range("E44").Value = '  -3.59%  '

$ws.Range("D45").Value2 = "'21.27"
$ws.Range("E45").Value = '  +1.20%  '

$ws.Range("E46").Value = '  +2.99%  '

$ws.Range("E47").Value = '  -2.60%  '

$ws.Range("D48").Value2 = "'0.0960"
$ws.Range("E48").Value = '  -1.93%  '

$ws.Range("D49").Value2 = "'0.0236"
$ws.Range("E49").Value = '  -1.36%  '

$ws.Range("E50").Value = '  -2.71%  '

$ws.Range("D51").Value2 = "'11.39"
$ws.Range("E51").Value = '  +0.00%  '
